# Add a new "PF/1.0.3" row to the meta-sheet, mirroring the existing
# per-environment layout (dev2/sit2/uat2/prod columns -> row 1 headers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.3"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"
